$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add two new rows of tlk-string data (row 62: "Armour of the Sten" header,
# row 63: its long description), following the exact same formatting pattern
# used by the previous similar pair of rows (58/59 - "TSM Alistairs Rose").
# ---------------------------------------------------------------------------

# Copy the direct cell formatting (fonts/fills/alignment) from the existing
# rows 58 (plain item-name row) and 59 (wrapped description row) so the new
# rows 62/63 end up with identical styles (s="4"/s="2" and s="4"/s="5").
$ws.Range("A58:B58").Copy()
$ws.Range("A62:B62").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A59:B59").Copy()
$ws.Range("A63:B63").PasteSpecial(-4122)   # xlPasteFormats

# Row 62 - "Armour of the Sten" item name.
$ws.Range("A62").Value = 6610060
$ws.Range("B62").Value = "Armour of the Sten"

# Row 63 - the long item description (wrapped, multi-line text).
$desc = "Traditionally worn by a Sten of the Berasaad, this armour is light and functional for the protection it offers.`r`nFinely linked chainmail, exotic hardened leather plates, and carefully ornamented metals show the craftsmanship that went into this piece.`r`nAs Qunari do not go by names, it is impossible to trace the origin of this armour, though it looks to have seen use in combat."
$ws.Range("A63").Value = 6610061
$ws.Range("B63").Value = $desc

# Match the source workbook's auto-fit row height for the wrapped text row.
$ws.Rows(63).RowHeight = 33.75

# ---------------------------------------------------------------------------
# Add a cell comment on A62 documenting the new entry, matching the style of
# the workbook's other "TSM ..." reviewer comments.
# ---------------------------------------------------------------------------
$excel.UserName = "Jim"
$ws.Range("A62").AddComment("Armour of the Sten")

# ---------------------------------------------------------------------------
# Update the active selection to match where the author left off editing.
# ---------------------------------------------------------------------------
$ws.Range("B70").Select()
